# Update cryptocurrency price (D) and 1h volume-change (E) columns
# with freshly scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (dropping significant trailing zeros, e.g. "6.20" -> 6.2) are pre-formatted
# as Text so the literal digit string from the source diff is preserved.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "59.306.86"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.576.71"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "555.96"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "141.82"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").Value = "2.584.34"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "6.73"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "0.166"
$ws.Range("E12").Value = "  +12.63%  "
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").Value = "3.027.80"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "59.296.37"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "22.95"
$ws.Range("E16").Value = "  +4.28%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "2.583.02"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").Value = "4.54"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "337.77"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").Value = "10.32"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "6.45"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("D24").Value = "0.477"
$ws.Range("E24").Value = "  +8.85%  "
$ws.Range("E25").Value = "  -5.10%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").Value = "7.38"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "0.0₃0774"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "6.20"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").Value = "159.01"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "19.06"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "4.09"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "0.894"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").Value = "37.44"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").Value = "0.851"
$ws.Range("E39").Value = "  -4.41%  "
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("D41").Value = "3.67"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").Value = "289.16"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").Value = "137.18"
$ws.Range("E43").Value = "  +7.22%  "
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "0.0973"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").Value = "10.68"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "0.0530"
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("D49").Value = "0.0234"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "18.65"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Value = "1.937.50"
$ws.Range("E51").Value = "  -1.16%  "

